$d = $word.ActiveDocument

$replacements = @(
    @("2024-04-13 Saturday", "2024-04-14 Sunday"),
    @("70×89=6230", "66×92=6072"),
    @("29×40=1160", "36×51=1836"),
    @("97×65=6305", "18×40=720"),
    @("84×13=1092", "91×31=2821"),
    @("98×51=4998", "53×16=848"),
    @("78×49=3822", "11×22=242"),
    @("69×21=1449", "49×27=1323"),
    @("75×62=4650", "72×86=6192"),
    @("90×80=7200", "59×57=3363"),
    @("57×29=1653", "71×40=2840"),
    @("45×70=3150", "20×26=520"),
    @("35×66=2310", "50×56=2800"),
    @("32×59=1888", "54×71=3834"),
    @("93×32=2976", "91×93=8463"),
    @("25×17=425", "96×89=8544"),
    @("34×60=2040", "44×87=3828"),
    @("70×94=6580", "38×32=1216"),
    @("12×45=540", "74×89=6586"),
    @("93×95=8835", "55×73=4015"),
    @("47×45=2115", "88×25=2200"),
    @("57×76=4332", "20×90=1800"),
    @("61×85=5185", "15×52=780"),
    @("29×13=377", "37×66=2442"),
    @("14×40=560", "74×31=2294"),
    @("95×53=5035", "97×99=9603")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Host "Done applying replacements"
